# Apply CRC (Common Rule / punctuation) and personal copy edits to the
# "definitions" lookup table on Sheet1.
#
# Each of these five edits only changes punctuation/quoting within an
# existing cell's text; Excel will naturally rebuild the shared-string
# table on save (moving the edited strings to the end of the table)
# exactly as shown in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31, column C - "Human subject" definition: semicolon -> colon
$ws.Range("C31").Value = "The Common Rule (45 CFR 46) definition of a human subject is a living individual about whom an investigator conducting research obtains: 1) Data through intervention or interaction with the individual, or 2) identifiable private information."

# Row 34, column C - "quasi-identifiers" definition: add comma after "gender)"
$ws.Range("C34").Value = "These variables do not alone identify a particular individual (e.g., ethnicity, gender), but, if combined with other information, they could be used to identify a participant"

# Row 38, column C - "Measure" definition: quote the word measure
$ws.Range("C38").Value = 'In this book, I use the term "measure" broadly to refer to a collection of items used to measure an outcome (e.g., an existing scale, an existing academic assessment).'

# Row 41, column C - "Normalize" definition: quote the word normalize (keep trailing space)
$ws.Range("C41").Value = 'In this book, the term "normalize" is used to refer to returning a value to its normal, or expected state '

# Row 43, column C - "Original data" definition: "First hand" -> "First-hand"
$ws.Range("C43").Value = "First-hand data that are generated/collected by the research team as part of the research study."

# Update the saved view state (scroll position + selection) to match the
# author's final cursor position in the workbook.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("F29").Select()
